$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.923.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.712.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4041"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4062"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.94"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.478"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08832"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.26"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.523"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.043"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001349"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.720.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.33"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07191"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.299"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.53"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.929.37"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.341"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.906"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.15"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.395"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +22.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.49"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "144.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.293"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.39%  "

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.939.97"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.29%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.273"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +14.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08808"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03186"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.293"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.034"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2883"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8488"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.89"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09498"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.23"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.476"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.74"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.736"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7480"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.83%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.21"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08418"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.62%  "

